# "Another Mini Fix before staging"
# The rate for the Navel product (row 2) was corrected from 140 to 400.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rates")

$ws.Range("B2").Value = 400
